# Sprint 1 Meeting 9 added
# Populates column K (meeting "3/11/2019: 6pm") with each team member's
# log entry for the new meeting, and adjusts the row heights that grew
# to fit the newly-added text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column K = new "Meeting 9" log entries (one per question/person row)
$ws.Range("K2").Value  = "Built application data validation for screens. Minor layout tweaks. Firebase updates."
$ws.Range("K3").Value  = "I worked on intents and implementation for the main menu, sign up page,e and about page."
$ws.Range("K4").Value  = "Finished the review page and help page"
$ws.Range("K5").Value  = "Successfully brought a Unity project over to Android studio and ran said project"
$ws.Range("K6").Value  = "Continue to look into firebase"
$ws.Range("K7").Value  = "I will work on my assigned issues. Every Issue assigned to me as well as the design and layout of every button and image in all pages for all issues of all assignee."
$ws.Range("K8").Value  = "Figure out mail gun for sending emails to the developers"
$ws.Range("K9").Value  = "Mapping character movement to UI buttons and integrating the Unity project into an existing Android studio project"
$ws.Range("K10").Value = "Difficulty reading/writing to firebase"
$ws.Range("K11").Value = "No, I am always working with all cylindrs pumping"
$ws.Range("K12").Value = "No"
$ws.Range("K13").Value = "Nothing is currently getting in the way of my work"
$ws.Range("K14").Value = "Establishing a connection to firebase"
$ws.Range("K15").Value = "Working together is better than alone!"
$ws.Range("K16").Value = "Team work is valuable and communication is key in order to save time"
$ws.Range("K17").Value = "Learned how to port a Unity project into Android studio"
$ws.Range("K18").Value = "Not currently"
$ws.Range("K19").Value = "No"
$ws.Range("K20").Value = "No"
$ws.Range("K21").Value = "No changes currently have to be made to the current plan for the project"

# Row heights grew on the rows whose new Meeting-9 text needed more
# vertical room than the tallest existing entry in that row.
$ws.Rows.Item(3).RowHeight  = 108
$ws.Rows.Item(5).RowHeight  = 97
$ws.Rows.Item(7).RowHeight  = 91.5
$ws.Rows.Item(9).RowHeight  = 81.5
$ws.Rows.Item(15).RowHeight = 73.5
$ws.Rows.Item(19).RowHeight = 157

# Move the view/selection to reflect where the author was working
$ws.Activate()
$ws.Range("J6").Select()
